# Updates the cryptocurrency price/volume table (generated from commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a "Price"-column (D) cell as literal TEXT, never
# letting Excel auto-convert a decimal-looking string ("314.35") into a number.
# Briefly flip the cell to Text format for the assignment, then restore the
# cell style to Normal so no visible formatting change is left behind.
function Set-PriceText($r, $c, $text) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row => @(new Price (D) value or $null, new Volume(1h) (E) value or $null)
$changes = @{
    2 = @("27.311.83", "  +0.88%  ")
    3 = @("1.824.05", "  -0.09%  ")
    4 = @($null, "  -0.14%  ")
    5 = @("314.35", "  +0.48%  ")
    6 = @($null, "  -0.10%  ")
    7 = @("0.4485", "  -1.80%  ")
    8 = @("0.3786", "  +1.69%  ")
    9 = @("0.07434", "  +1.78%  ")
    10 = @("0.8861", "  +3.04%  ")
    11 = @("20.97", "  +0.41%  ")
    12 = @("1.820.50", "  -0.21%  ")
    13 = @("6.731", $null)
    14 = @("5.451", "  +1.82%  ")
    15 = @("93.54", "  +0.54%  ")
    16 = @("0.07130", "  +0.34%  ")
    17 = @("1.001", "  -0.27%  ")
    18 = @("0.000008808", "  -0.38%  ")
    19 = @($null, "  -0.09%  ")
    20 = @("15.14", "  +0.89%  ")
    21 = @("27.333.36", "  +0.97%  ")
    22 = @("5.389", "  +3.86%  ")
    23 = @("10.96", "  -0.03%  ")
    24 = @($null, "  -1.77%  ")
    25 = @("151.68", "  -0.06%  ")
    26 = @("2.312", "  +4.04%  ")
    27 = @("18.64", "  +0.98%  ")
    28 = @("5.388", "  +2.20%  ")
    29 = @("117.81", "  +0.34%  ")
    30 = @("0.08899", "  +0.19%  ")
    31 = @("0.7923", "  +4.54%  ")
    32 = @("1.202", "  +0.77%  ")
    33 = @("4.601", "  +2.91%  ")
    34 = @("2.916", "  -1.59%  ")
    35 = @($null, "  -0.16%  ")
    36 = @("1.113", "  +1.06%  ")
    37 = @($null, "  +0.57%  ")
    38 = @("0.05303", "  +0.39%  ")
    39 = @("7.353", "  +2.29%  ")
    40 = @("0.5339", "  -0.29%  ")
    41 = @("2.869", "  -0.50%  ")
    44 = @("8.671", "  +1.12%  ")
    45 = @("0.5068", "  -3.35%  ")
    46 = @("10.67", "  -0.51%  ")
    47 = @("1.700", "  +1.35%  ")
    48 = @("105.30", "  -0.33%  ")
    49 = @("0.9999", "  -0.14%  ")
    50 = @($null, "  -0.08%  ")
    51 = @("66.00", "  +3.94%  ")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    if ($null -ne $vals[0]) { Set-PriceText $row 4 $vals[0] }
    if ($null -ne $vals[1]) { $ws.Cells.Item($row, 5).Value = $vals[1] }
}

# Rows 42 & 43 swap ranking positions (Algorand <-> RenderToken); rank in column A
# is unchanged, only Coin/Link/Price/Volume(1h) move between the two rows.
$ws.Cells.Item(42, 2).Value = "RenderToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-PriceText 42 4 "2.338"
$ws.Cells.Item(42, 5).Value = "  +18.57%  "

$ws.Cells.Item(43, 2).Value = "Algorand"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-PriceText 43 4 "0.1716"
$ws.Cells.Item(43, 5).Value = "  -0.11%  "

